# Append daily COVID-19 case data through 2021-06-28 (Sassuolo report).
# New rows 270-301 extend columns A (date), B (nuovi pos.),
# C (somma mobile 7gg.) and D (somma mobile 7gg. per 100mila abitanti).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries the date-time cell style (s="2"); extend it down
# across the new rows by copying the format from the last existing row
# before the values themselves are written (Copy then overwrite Value).
$ws.Range("A269").Copy($ws.Range("A270:A301"))

# r, A (date serial), B (nuovi pos.), C (somma mobile 7gg.), D (per 100k ab.)
$newData = @(
    @(270, 44344, 3, 11, 27.31019415065296),
    @(271, 44345, 5, 13, 32.27568399622623),
    @(272, 44346, 3, 14, 34.75842891901286),
    @(273, 44347, 3, 16, 39.72391876458613),
    @(274, 44348, 0, 16, 39.72391876458613),
    @(275, 44349, 1, 16, 39.72391876458613),
    @(276, 44350, 2, 17, 42.20666368737276),
    @(277, 44351, 3, 17, 42.20666368737276),
    @(278, 44352, 4, 16, 39.72391876458613),
    @(279, 44353, 0, 13, 32.27568399622623),
    @(280, 44354, 1, 11, 27.31019415065296),
    @(281, 44355, 1, 12, 29.79293907343959),
    @(282, 44356, 2, 13, 32.27568399622623),
    @(283, 44357, 1, 12, 29.79293907343959),
    @(284, 44358, 0, 9, 22.34470430507969),
    @(285, 44359, 2, 7, 17.37921445950643),
    @(286, 44360, 0, 7, 17.37921445950643),
    @(287, 44361, 0, 6, 14.8964695367198),
    @(288, 44362, 1, 6, 14.8964695367198),
    @(289, 44363, 0, 4, 9.930979691146533),
    @(290, 44364, 0, 3, 7.448234768359899),
    @(291, 44365, 0, 3, 7.448234768359899),
    @(292, 44366, 0, 1, 2.482744922786633),
    @(293, 44367, 0, 1, 2.482744922786633),
    @(294, 44368, 0, 1, 2.482744922786633),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 1, 1, 2.482744922786633),
    @(297, 44371, 1, 2, 4.965489845573266),
    @(298, 44372, 0, 2, 4.965489845573266),
    @(299, 44373, 1, 3, 7.448234768359899),
    @(300, 44374, 1, 4, 9.930979691146533),
    @(301, 44375, 0, 4, 9.930979691146533),
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
